# Applies the "Updated symbol list" data refresh (Fri Jan 20 07:49:54 UTC 2023)
# to the cryptos worksheet: Price (column D) and Volume(1h) (column E) cells
# are refreshed with the newer scrape values for the affected rows.
#
# The cells hold plain text values (e.g. "288.42", "1.02%"), not numbers, so
# each cell's number format is reset to Text ("@") before the new value is
# written. That keeps Excel from auto-coercing a numeric- or percent-looking
# string into a real number/percentage -- the stored cell content stays an
# exact text string, matching the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "288.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.02%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.934"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.67%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07354"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.65%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.263"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "25.65%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.727"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.55%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.729"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9034"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.83%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08980"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "16.85%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1685"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08180"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03109"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.16%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09951"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.71%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.34%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005732"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.24%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.489"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.56%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.061"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.04%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.42%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.99%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.154"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.94%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2110"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-9.14%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04554"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.02%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001209"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.65%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004162"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.68%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.94%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01581"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.42%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04454"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.11%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007390"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.24%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009593"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "24.86%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1327"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.48%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002220"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "16.78%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008502"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-7.70%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.46%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.100"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-6.54%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002000"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-33.36%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
